$d = $word.ActiveDocument
$d.Content.Find.Execute("September", $true, $false, $false, $false, $false, $true, 1, $false, "Enero", 2)
$d.Content.Find.Execute("2016", $true, $false, $false, $false, $false, $true, 1, $false, "2017", 2)
